# Weekly update: insert two new daily price rows for "Ají" into the
# Macroferia Regional de Talca sheet.
#
# Row 100 (new) is inserted above the existing data (pushing the former
# rows 100-160 down to 101-161), and a second new row is inserted at the
# (post-shift) position 156, pushing the remaining former rows down once
# more (final sheet ends at row 162, matching dimension A1:R162).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the first new row, above the current row 100 ---------------
$ws.Rows.Item(100).Insert()

$ws.Range("A100").Value = 5
$ws.Range("B100").Value = "Macroferia Regional de Talca"
$ws.Range("C100").Value = "Maule"
$ws.Range("D100").Value = 44567
$ws.Range("E100").Value = 7
$ws.Range("F100").Value = 100112021
$ws.Range("G100").Value = "Ají"
$ws.Range("H100").Value = "Americana (o)"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 150
$ws.Range("K100").Value = 15000
$ws.Range("L100").Value = 15000
$ws.Range("M100").Value = 15000
$ws.Range("N100").Value = "$/caja 14 kilos"
$ws.Range("O100").Value = "Región del Maule"
$ws.Range("P100").Value = 1071
$ws.Range("Q100").Value = 14
$ws.Range("R100").Value = "Hortaliza"

# --- Insert the second new row, at (post-shift) row 156 -----------------
$ws.Rows.Item(156).Insert()

$ws.Range("A156").Value = 5
$ws.Range("B156").Value = "Macroferia Regional de Talca"
$ws.Range("C156").Value = "Maule"
$ws.Range("D156").Value = 44568
$ws.Range("E156").Value = 7
$ws.Range("F156").Value = 100112021
$ws.Range("G156").Value = "Ají"
$ws.Range("H156").Value = "Americana (o)"
$ws.Range("I156").Value = "Primera"
$ws.Range("J156").Value = 100
$ws.Range("K156").Value = 15000
$ws.Range("L156").Value = 15000
$ws.Range("M156").Value = 15000
$ws.Range("N156").Value = "$/caja 14 kilos"
$ws.Range("O156").Value = "Región del Maule"
$ws.Range("P156").Value = 1071
$ws.Range("Q156").Value = 14
$ws.Range("R156").Value = "Hortaliza"
